# TP3E4_PROJ.xlsx - "Avances en el diseno"
#
# - Fill in the component values that were chosen for R3/R4 and the
#   op-amp question for U1/U2 on the TRABAJO sheet.
# - Make TRABAJO the active sheet/tab (it was MPX2010DP before), leaving
#   the selection on C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRABAJO")

# Switch to the TRABAJO sheet (was on MPX2010DP before the edit).
$ws.Activate()

# New values entered in column C for the R3/R4/U1/U2 subtasks.
$ws.Range("C7").Value = 12
$ws.Range("C8").Value = 1522
$ws.Range("C9").Value = "LM833 con +VCC=15 y -VCC=0?"
$ws.Range("C10").Value = "LM833 con +VCC=15 y -VCC=0?"

# Leave the selection where it ended up after entering the data.
$ws.Range("C6").Select()
